# Add a new week's worth of Plátano price data to the "Feria Lagunitas de
# Puerto Montt" sheet. This inserts two new rows right before the current
# row 891 (pushing all the following rows down by two, from 891-933 to
# 893-935) and fills the two new rows with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 891; this shifts every row currently at
# 891-933 down to 893-935, carrying all of their existing data and
# formatting with them (matching the target workbook exactly).
$ws.Rows("891:892").Insert()

# New row 891: "Pintón" quality entry for the new week (2023-08-09).
$row891 = New-Object 'object[,]' 1,20
$row891[0,0]  = 4
$row891[0,1]  = "Feria Lagunitas de Puerto Montt"
$row891[0,2]  = "Los Lagos"
$row891[0,3]  = 45147
$row891[0,4]  = 10
$row891[0,5]  = "Fruta"
$row891[0,6]  = 100108
$row891[0,7]  = "Tropicales y subtropicales"
$row891[0,8]  = 100108006
$row891[0,9]  = "Plátano"
$row891[0,10] = "Sin especificar"
$row891[0,11] = "Pintón"
$row891[0,12] = 400
$row891[0,13] = 17000
$row891[0,14] = 17000
$row891[0,15] = 17000
$row891[0,16] = "`$/caja 20 kilos"
$row891[0,17] = "Ecuador"
$row891[0,18] = 850
$row891[0,19] = 20
$ws.Range("A891:T891").Value = $row891

# New row 892: "Primera Pintón" quality entry for the same new week.
$row892 = New-Object 'object[,]' 1,20
$row892[0,0]  = 4
$row892[0,1]  = "Feria Lagunitas de Puerto Montt"
$row892[0,2]  = "Los Lagos"
$row892[0,3]  = 45147
$row892[0,4]  = 10
$row892[0,5]  = "Fruta"
$row892[0,6]  = 100108
$row892[0,7]  = "Tropicales y subtropicales"
$row892[0,8]  = 100108006
$row892[0,9]  = "Plátano"
$row892[0,10] = "Sin especificar"
$row892[0,11] = "Primera Pintón"
$row892[0,12] = 400
$row892[0,13] = 20000
$row892[0,14] = 20000
$row892[0,15] = 20000
$row892[0,16] = "`$/caja 20 kilos"
$row892[0,17] = "Ecuador"
$row892[0,18] = 1000
$row892[0,19] = 20
$ws.Range("A892:T892").Value = $row892
